$wb = $excel.ActiveWorkbook

# --- Update data on "reatores" sheet ---
$reatores = $wb.Worksheets.Item("reatores")

$reatores.Range("C2").Value = 20
$reatores.Range("D2").Value = 37.5
$reatores.Range("E2").Value = 58

$reatores.Range("J2").Formula = "=escolas!C2"
$reatores.Range("J3").Formula = "=J2"

# --- Update selection on "escolas" sheet (now not the active tab) ---
$escolas = $wb.Worksheets.Item("escolas")
$escolas.Activate() | Out-Null
$escolas.Range("A1:I2").Select() | Out-Null

# --- Activate "reatores" sheet and set its selection (becomes the active tab) ---
$reatores.Activate() | Out-Null
$reatores.Range("E3").Select() | Out-Null
